# VITA 57.4 LPC ADC/DAC Pinout - relabel the GPIO / DAC_DB / DAC_CLK cells
# on the "Table 5-1 LPC ADC|DAC Pinout" sheet, and update the saved
# selection/active-cell on both sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws  = $wb.Worksheets.Item(2)

# --- Relabel cells: GPIO numbering shifts by +2 (room made for GPIO0/GPIO1
#     that used to be the DAC_CLK_P/N clock pins), and the DAC_DB/DAC_WRT
#     pins shuffle down to make room, with DAC_CLK_P/N moving to where
#     DAC_WRTB/DAC_WRTA used to be. ---

$ws.Range("E3").Interior.Color = 10092543   # FFFFFF99 (same fill as other GPIO cells)
$ws.Range("E3").Value = "GPIO0"

$ws.Range("E4").Interior.Color = 10092543   # FFFFFF99
$ws.Range("E4").Value = "GPIO1"

$ws.Range("H9").Value = "GPIO2"
$ws.Range("H10").Value = "GPIO3"
$ws.Range("I11").Value = "GPIO4"
$ws.Range("H12").Value = "GPIO6"
$ws.Range("I12").Value = "GPIO5"
$ws.Range("H13").Value = "GPIO7"
$ws.Range("H15").Value = "GPIO8"
$ws.Range("I15").Value = "GPIO10"
$ws.Range("H16").Value = "GPIO9"
$ws.Range("I16").Value = "GPIO11"
$ws.Range("H18").Value = "GPIO12"
$ws.Range("H19").Value = "GPIO13"
$ws.Range("I19").Value = "GPIO14"
$ws.Range("I20").Value = "GPIO15"
$ws.Range("H21").Value = "GPIO16"
$ws.Range("H22").Value = "GPIO17"
$ws.Range("I23").Value = "GPIO18"
$ws.Range("H24").Value = "DAC_DB0"
$ws.Range("I24").Value = "GPIO19"
$ws.Range("E25").Value = "DAC_DB2"
$ws.Range("H25").Value = "DAC_DB1"
$ws.Range("D26").Value = "DAC_DB4"
$ws.Range("E26").Value = "DAC_DB3"
$ws.Range("D27").Value = "DAC_DB5"

$ws.Range("H27").Interior.Color = 52479     # FFFFCC00 (same fill as ADC_CLK cells)
$ws.Range("H27").Value = "DAC_CLK_P (100MHz)"

$ws.Range("E28").Value = "DAC_DB6"

$ws.Range("H28").Interior.Color = 52479     # FFFFCC00
$ws.Range("H28").Value = "DAC_CLK_N (100MHz)"

$ws.Range("D29").Value = "DAC_DB8"
$ws.Range("E29").Value = "DAC_DB7"
$ws.Range("D30").Value = "DAC_DB9"
$ws.Range("E31").Value = "DAC_WRTB"
$ws.Range("E32").Value = "DAC_WRTA"

# --- Update view state: scroll position + active cell on each sheet ---

[void]$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws1.Range("D30").Select()

[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("N32").Select()
